$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.405.04'
$ws.Range("E2").Value = '  -4.18%  '

$ws.Range("D3").Value = '2.533.87'
$ws.Range("E3").Value = '  -3.27%  '

$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '505.92'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -4.29%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.42'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -7.46%  '

$ws.Range("E7").Value = '  +0.04%  '

$ws.Range("E8").Value = '  -5.13%  '

$ws.Range("D9").Value = '2.536.12'
$ws.Range("E9").Value = '  -3.33%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.17'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -7.57%  '

$ws.Range("E11").Value = '  -7.01%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.331'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -4.81%  '

$ws.Range("E13").Value = '  -0.59%  '

$ws.Range("D14").Value = '2.977.26'
$ws.Range("E14").Value = '  -3.27%  '

$ws.Range("D15").Value = '58.378.64'
$ws.Range("E15").Value = '  -4.20%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '20.64'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -4.89%  '

$ws.Range("E17").Value = '  -6.26%  '

$ws.Range("D18").Value = '2.537.44'
$ws.Range("E18").Value = '  -3.22%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.53'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -5.31%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '335.47'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -5.34%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.06'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -5.18%  '

$ws.Range("E22").Value = '  -0.08%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.94'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -4.74%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '60.04'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -2.33%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.407'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -5.26%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +0.07%  '

$ws.Range("E27").Value = '  -5.18%  '

$ws.Range("D28").Value = '2.647.15'
$ws.Range("E28").Value = '  -3.22%  '

$ws.Range("D29").Value = '0.0₃0786'
$ws.Range("E29").Value = '  -9.21%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.93'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -6.13%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.999'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -0.05%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '149.79'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -0.66%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.85'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -5.51%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '18.49'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -4.98%  '

$ws.Range("E35").Value = '  -5.22%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.937'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +5.00%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.91'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -6.24%  '

$ws.Range("E38").Value = '  -7.68%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '36.04'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -1.07%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.824'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -11.56%  '

$ws.Range("E41").Value = '  -6.74%  '

$ws.Range("B42").Value = 'Bittensor'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '283.32'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -4.03%  '

$ws.Range("B43").Value = 'Filecoin'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.52'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -7.29%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0997'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -3.13%  '

$ws.Range("E45").Value = '  +0.00%  '

$ws.Range("E46").Value = '  -5.93%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0534'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -4.92%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '18.65'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -5.48%  '

$ws.Range("E49").Value = '  -0.49%  '

$ws.Range("E50").Value = '  -5.26%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.52'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -7.78%  '
